$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1,4).Value = "Euro"
$ws.Cells.Item(2,4).Value = 141.73
$ws.Cells.Item(3,4).Value = 38.56
$ws.Cells.Item(4,4).Value = 40
$ws.Cells.Item(5,4).Value = 55.4
$ws.Cells.Item(6,4).Value = 50.0
$ws.Cells.Item(7,4).Value = 823
$ws.Cells.Item(8,4).Value = 630.0
$ws.Cells.Item(9,4).Value = 150.36
$ws.Cells.Item(10,4).Value = 100.32
$ws.Cells.Item(11,4).Value = 30.0
$ws.Cells.Item(12,4).Value = 130.0
$ws.Cells.Item(13,4).Value = 1185.57
$ws.Cells.Item(14,4).Value = 1070.56
$ws.Cells.Item(15,4).Value = 4.5
$ws.Cells.Item(16,4).Value = 135.26
$ws.Cells.Item(17,4).Value = 13.75
$ws.Cells.Item(18,4).Value = 5.5
$ws.Cells.Item(19,4).Value = 113.4
$ws.Cells.Item(20,4).Value = 116.0
$ws.Cells.Item(21,4).Value = 71.95
$ws.Cells.Item(22,4).Value = 10.0
$ws.Cells.Item(23,4).Value = 243.56
$ws.Cells.Item(24,4).Value = 24.0
$ws.Cells.Item(25,4).Value = 542.56
$ws.Cells.Item(26,4).Value = 533.0
$ws.Cells.Item(27,4).Value = 344.36
$ws.Cells.Item(28,4).Value = 255.16
$ws.Cells.Item(29,4).Value = 91
$ws.Cells.Item(30,4).Value = 54.9
$ws.Cells.Item(31,4).Value = 7004.4
$ws.Columns.Item(5).Delete()

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1,4).Value = "Euro"
$ws.Cells.Item(2,4).Value = 544
$ws.Cells.Item(3,4).Value = 1060.05
$ws.Cells.Item(4,4).Value = 1028.81
$ws.Cells.Item(5,4).Value = 52.36
$ws.Cells.Item(6,4).Value = 50.0
$ws.Cells.Item(7,4).Value = 2760
$ws.Cells.Item(8,4).Value = 256
$ws.Cells.Item(9,4).Value = 635.0
$ws.Cells.Item(10,4).Value = 635.0
$ws.Cells.Item(11,4).Value = 635.0
$ws.Cells.Item(12,4).Value = 635.0
$ws.Cells.Item(13,4).Value = 226.89
$ws.Cells.Item(14,4).Value = 6.2
$ws.Cells.Item(15,4).Value = 600
$ws.Cells.Item(16,4).Value = 530
$ws.Cells.Item(17,4).Value = 60
$ws.Cells.Item(18,4).Value = 920
$ws.Cells.Item(19,4).Value = 20
$ws.Cells.Item(20,4).Value = 326.5
$ws.Cells.Item(21,4).Value = 10980.81
$ws.Columns.Item(5).Delete()

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1,4).Value = "Euro"
$ws.Cells.Item(26,4).Value = "Euro"
$ws.Cells.Item(2,4).Value = 544
$ws.Cells.Item(3,4).Value = 1060.05
$ws.Cells.Item(4,4).Value = 1028.81
$ws.Cells.Item(5,4).Value = 52.36
$ws.Cells.Item(6,4).Value = 50.0
$ws.Cells.Item(7,4).Value = 2760
$ws.Cells.Item(8,4).Value = 256
$ws.Cells.Item(9,4).Value = 635.0
$ws.Cells.Item(10,4).Value = 635.0
$ws.Cells.Item(11,4).Value = 635.0
$ws.Cells.Item(12,4).Value = 635.0
$ws.Cells.Item(13,4).Value = 226.89
$ws.Cells.Item(14,4).Value = 6.2
$ws.Cells.Item(15,4).Value = 600
$ws.Cells.Item(16,4).Value = 530
$ws.Cells.Item(17,4).Value = 60
$ws.Cells.Item(18,4).Value = 920
$ws.Cells.Item(19,4).Value = 20
$ws.Cells.Item(20,4).Value = 326.5
$ws.Cells.Item(21,4).Value = 10980.81
$ws.Cells.Item(27,4).Value = 141.73
$ws.Cells.Item(28,4).Value = 38.56
$ws.Cells.Item(29,4).Value = 40
$ws.Cells.Item(30,4).Value = 55.4
$ws.Cells.Item(31,4).Value = 50.0
$ws.Cells.Item(32,4).Value = 823
$ws.Cells.Item(33,4).Value = 630.0
$ws.Cells.Item(34,4).Value = 150.36
$ws.Cells.Item(35,4).Value = 100.32
$ws.Cells.Item(36,4).Value = 30.0
$ws.Cells.Item(37,4).Value = 130.0
$ws.Cells.Item(38,4).Value = 1185.57
$ws.Cells.Item(39,4).Value = 1070.56
$ws.Cells.Item(40,4).Value = 4.5
$ws.Cells.Item(41,4).Value = 135.26
$ws.Cells.Item(42,4).Value = 13.75
$ws.Cells.Item(43,4).Value = 5.5
$ws.Cells.Item(44,4).Value = 113.4
$ws.Cells.Item(45,4).Value = 116.0
$ws.Cells.Item(46,4).Value = 71.95
$ws.Cells.Item(47,4).Value = 10.0
$ws.Cells.Item(48,4).Value = 243.56
$ws.Cells.Item(49,4).Value = 24.0
$ws.Cells.Item(50,4).Value = 542.56
$ws.Cells.Item(51,4).Value = 533.0
$ws.Cells.Item(52,4).Value = 344.36
$ws.Cells.Item(53,4).Value = 255.16
$ws.Cells.Item(54,4).Value = 91
$ws.Cells.Item(55,4).Value = 54.9
$ws.Cells.Item(56,4).Value = 7004.4
$ws.Columns.Item(5).Delete()
